$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 59, shifting existing rows 59-60 down to 60-61.
$ws.Rows.Item(59).Insert()

# New row 59 data (latest week), copying formatting/values from the row that
# used to be at 59 (now 60) where unchanged, and updating the rest per diff.
$ws.Cells.Item(59, 1).Value = 5
$ws.Cells.Item(59, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(59, 3).Value = "Maule"
$ws.Cells.Item(59, 4).Value = 44595
$ws.Cells.Item(59, 5).Value = 7
$ws.Cells.Item(59, 6).Value = "Fruta"
$ws.Cells.Item(59, 7).Value = 100101
$ws.Cells.Item(59, 8).Value = "Berries"
$ws.Cells.Item(59, 9).Value = 100101001
$ws.Cells.Item(59, 10).Value = "Arándano (blue)"
$ws.Cells.Item(59, 11).Value = "Sin especificar"
$ws.Cells.Item(59, 12).Value = "Primera"
$ws.Cells.Item(59, 13).Value = 150
$ws.Cells.Item(59, 14).Value = 3500
$ws.Cells.Item(59, 15).Value = 3500
$ws.Cells.Item(59, 16).Value = 3500
$ws.Cells.Item(59, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(59, 18).Value = "Provincia de Linares"
$ws.Cells.Item(59, 19).Value = 1750
$ws.Cells.Item(59, 20).Value = 2
